$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibition) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 838
$ws1.Cells.Item(3, 6).Value = 13684
$ws1.Cells.Item(4, 6).Value = 13473
$ws1.Cells.Item(5, 6).Value = 1045
$ws1.Cells.Item(6, 6).Value = 801
$ws1.Cells.Item(8, 6).Value = 587
$ws1.Cells.Item(9, 6).Value = 78
$ws1.Cells.Item(11, 6).Value = 43
$ws1.Cells.Item(12, 6).Value = 746
$ws1.Cells.Item(15, 6).Value = 81
$ws1.Cells.Item(16, 6).Value = 68
$ws1.Cells.Item(17, 6).Value = 110
$ws1.Cells.Item(19, 6).Value = 505
$ws1.Cells.Item(20, 6).Value = 422
$ws1.Cells.Item(21, 6).Value = 374
$ws1.Cells.Item(22, 6).Value = 306
$ws1.Cells.Item(23, 6).Value = 4
$ws1.Cells.Item(24, 6).Value = 819
$ws1.Cells.Item(25, 6).Value = 70

# Append new row 26 (2024-08-17 Dragon Ball exhibit) to 展览
$ws1.Cells.Item(25, 1).Copy($ws1.Cells.Item(26, 1))
$ws1.Cells.Item(26, 1).Value = 25
$ws1.Cells.Item(26, 2).NumberFormat = "@"
$ws1.Cells.Item(26, 2).Value = "2024-08-17"
$ws1.Cells.Item(26, 2).Style = "Normal"
$ws1.Cells.Item(26, 3).Value = "广州·鸟山明作品《龙珠》40周年only纪念展"
$ws1.Cells.Item(26, 4).Value = "逸景路462号珠江国际纺织城d区6层 珠江时尚馆"
$ws1.Cells.Item(26, 5).Value = "2024.08.17 10:00-08.17 17:30"
$ws1.Cells.Item(26, 6).Value = 0
$ws1.Cells.Item(26, 7).Value = 68
$ws1.Cells.Item(26, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86780"
$ws1.Cells.Item(26, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/4k7Thger1717147185584.jpeg"

# ---- Sheet: 演出 (Performance) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(3, 6).Value = 22
$ws2.Cells.Item(5, 6).Value = 68
$ws2.Cells.Item(7, 6).Value = 1420

# ---- Sheet: 全部类型 (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 838
$ws4.Cells.Item(4, 6).Value = 13684
$ws4.Cells.Item(5, 6).Value = 13473
$ws4.Cells.Item(6, 6).Value = 1045
$ws4.Cells.Item(7, 6).Value = 801
$ws4.Cells.Item(9, 6).Value = 587
$ws4.Cells.Item(10, 6).Value = 78
$ws4.Cells.Item(12, 6).Value = 43
$ws4.Cells.Item(13, 6).Value = 746
$ws4.Cells.Item(15, 6).Value = 22
$ws4.Cells.Item(18, 6).Value = 81
$ws4.Cells.Item(19, 6).Value = 68
$ws4.Cells.Item(20, 6).Value = 110
$ws4.Cells.Item(23, 6).Value = 68
$ws4.Cells.Item(26, 6).Value = 505
$ws4.Cells.Item(27, 6).Value = 422
$ws4.Cells.Item(28, 6).Value = 374
$ws4.Cells.Item(29, 6).Value = 306
$ws4.Cells.Item(30, 6).Value = 4
$ws4.Cells.Item(31, 6).Value = 819
$ws4.Cells.Item(33, 6).Value = 1420
$ws4.Cells.Item(36, 6).Value = 71

# Prime row 39/40 A-column formatting (bold/border/center) from an existing styled cell
# BEFORE overwriting any values, so the source cell (A38) is untouched at copy time.
$ws4.Cells.Item(38, 1).Copy($ws4.Cells.Item(40, 1))
$ws4.Cells.Item(38, 1).Copy($ws4.Cells.Item(39, 1))

# Shift old row 39 (孟京辉 "一个陌生女人的来信") down to row 40
$ws4.Cells.Item(40, 1).Value = 39
$ws4.Cells.Item(40, 2).NumberFormat = "@"
$ws4.Cells.Item(40, 2).Value = "2024-08-30"
$ws4.Cells.Item(40, 2).Style = "Normal"
$ws4.Cells.Item(40, 3).Value = "广州·孟京辉经典戏剧作品·黄湘丽主演《一个陌生女人的来信》"
$ws4.Cells.Item(40, 4).Value = "广州市越秀区人民北路696号 广州友谊剧院"
$ws4.Cells.Item(40, 5).Value = "2024.08.30 19:30-08.31 16:30"
$ws4.Cells.Item(40, 6).Value = 7
$ws4.Cells.Item(40, 7).Value = 100
$ws4.Cells.Item(40, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84570"
$ws4.Cells.Item(40, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/SscDFm1z1713177818070.jpeg"

# Write the new row 39 (2024-08-17 Dragon Ball exhibit) into 全部类型
$ws4.Cells.Item(39, 1).Value = 38
$ws4.Cells.Item(39, 2).NumberFormat = "@"
$ws4.Cells.Item(39, 2).Value = "2024-08-17"
$ws4.Cells.Item(39, 2).Style = "Normal"
$ws4.Cells.Item(39, 3).Value = "广州·鸟山明作品《龙珠》40周年only纪念展"
$ws4.Cells.Item(39, 4).Value = "逸景路462号珠江国际纺织城d区6层 珠江时尚馆"
$ws4.Cells.Item(39, 5).Value = "2024.08.17 10:00-08.17 17:30"
$ws4.Cells.Item(39, 6).Value = 0
$ws4.Cells.Item(39, 7).Value = 68
$ws4.Cells.Item(39, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86780"
$ws4.Cells.Item(39, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/4k7Thger1717147185584.jpeg"

Write-Host "edit applied"
